$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.004754311418771289
$ws.Range("C2").Value = 0.005389076851612592
$ws.Range("D2").Value = 0.00585536679281113
$ws.Range("E2").Value = 0.005670665460181294
$ws.Range("B3").Value = 2.95150647669481
$ws.Range("C3").Value = 4.754240764293362
$ws.Range("D3").Value = 5.585856170897704
$ws.Range("E3").Value = 5.992451452608591
$ws.Range("B4").Value = -0.01243754993729038
$ws.Range("C4").Value = -0.0112630984422302
$ws.Range("D4").Value = -0.01362021297418346
$ws.Range("E4").ClearContents()
$ws.Range("B5").Value = -5.488083401140036
$ws.Range("C5").Value = -5.696757530196266
$ws.Range("D5").Value = -6.256129489798984
$ws.Range("E5").ClearContents()
$ws.Range("B6").Value = 0.005116256475096994
$ws.Range("C6").Value = 0.001959664171740027
$ws.Range("D6").Value = 0.001743640555890115
$ws.Range("E6").Value = 0.001913477689023771
$ws.Range("B7").Value = 2.005235405987511
$ws.Range("C7").Value = 1.002575273488103
$ws.Range("D7").Value = 0.9461598357465517
$ws.Range("E7").Value = 1.08110823842823
$ws.Range("B8").Value = 0.006451957805491795
$ws.Range("C8").Value = 0.006356541366461522
$ws.Range("D8").Value = 0.006696391702962223
$ws.Range("E8").Value = 0.005730842260671652
$ws.Range("B9").Value = 3.858141706739029
$ws.Range("C9").Value = 5.400859541094776
$ws.Range("D9").Value = 6.478674860957485
$ws.Range("E9").Value = 6.379229698710517
$ws.Range("B10").Value = -0.01036271929057376
$ws.Range("C10").Value = -0.0109826750875756
$ws.Range("D10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("B11").Value = -4.527308866029893
$ws.Range("C11").Value = -5.376324566222857
$ws.Range("D11").ClearContents()
$ws.Range("E11").ClearContents()
$ws.Range("B12").Value = 0.001462100388132495
$ws.Range("C12").Value = 0.0006355493837071288
$ws.Range("D12").Value = 0.0006178379737913425
$ws.Range("E12").Value = 0.002467234319116947
$ws.Range("B13").Value = 0.5789676885702968
$ws.Range("C13").Value = 0.3179093725469145
$ws.Range("D13").Value = 0.3403269381641657
$ws.Range("E13").Value = 1.389384295598894
$ws.Range("B14").Value = 0.006695666715188392
$ws.Range("C14").Value = 0.006681371175403872
$ws.Range("D14").Value = 0.006641293356938142
$ws.Range("E14").Value = 0.00569243078486502
$ws.Range("B15").Value = 3.744964528610712
$ws.Range("C15").Value = 5.432774179440074
$ws.Range("D15").Value = 6.37000633007165
$ws.Range("E15").Value = 6.026796462879193
$ws.Range("B16").Value = -0.01030321078678046
$ws.Range("C16").Value = -0.01047400383540542
$ws.Range("D16").ClearContents()
$ws.Range("E16").ClearContents()
$ws.Range("B17").Value = -4.315604500821167
$ws.Range("C17").Value = -4.837054944352429
$ws.Range("D17").ClearContents()
$ws.Range("E17").ClearContents()
$ws.Range("B18").Value = 0.0008916621096052218
$ws.Range("C18").Value = -0.0002745196274605644
$ws.Range("D18").Value = 0.0003834457850681817
$ws.Range("E18").Value = 0.002107137268187614
$ws.Range("B19").Value = 0.3252430865880111
$ws.Range("C19").Value = -0.1406477203802823
$ws.Range("D19").Value = 0.2173789790247193
$ws.Range("E19").Value = 1.14675825605596
$ws.Range("B20").Value = 0.006366746531862169
$ws.Range("C20").Value = 0.006750748498099018
$ws.Range("D20").Value = 0.006539269801902266
$ws.Range("E20").Value = 0.005772245651577344
$ws.Range("B21").Value = 3.617560667053156
$ws.Range("C21").Value = 5.44956356670746
$ws.Range("D21").Value = 5.984021443800612
$ws.Range("E21").Value = 5.817743972605122
$ws.Range("B22").Value = -0.0105329681839226
$ws.Range("C22").Value = -0.01171470416059327
$ws.Range("D22").ClearContents()
$ws.Range("E22").ClearContents()
$ws.Range("B23").Value = -4.226845728165167
$ws.Range("C23").Value = -4.569805526302098
$ws.Range("D23").ClearContents()
$ws.Range("E23").ClearContents()
$ws.Range("B24").Value = 0.00142587876152263
$ws.Range("C24").Value = 0.0003046138437965086
$ws.Range("D24").Value = 0.000639535431532716
$ws.Range("E24").Value = 0.001518637830056985
$ws.Range("B25").Value = 0.5319866226561412
$ws.Range("C25").Value = 0.151309745251501
$ws.Range("D25").Value = 0.3208162709328659
$ws.Range("E25").Value = 0.7651011408357508
